$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Label" header in column H, row 1, matching header style of existing headers
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Updated D/E (refit) values and new H (Label: 0 = Control, 1 = MDD) column values
# Iterations = 100 block (rows 2-11)
$ws.Range("D2").Value = 0.4879964957685827
$ws.Range("E2").Value = 0.4879964957685827
$ws.Range("H2").Value = 0

$ws.Range("D3").Value = 0.5154319343229138
$ws.Range("E3").Value = 0.5154319343229138
$ws.Range("H3").Value = 0

$ws.Range("D4").Value = 0.1309795708441029
$ws.Range("E4").Value = 0.1309795708441029
$ws.Range("H4").Value = 0

$ws.Range("D5").Value = 0.2058723500158631
$ws.Range("E5").Value = 0.2058723500158631
$ws.Range("H5").Value = 0

$ws.Range("D6").Value = 0.01063413504050413
$ws.Range("E6").Value = 0.01063413504050413
$ws.Range("H6").Value = 0

$ws.Range("D7").Value = 0.4856912080921306
$ws.Range("E7").Value = 0.5143087919078694
$ws.Range("H7").Value = 1

$ws.Range("H8").Value = 1

$ws.Range("D9").Value = 0.3769155917912217
$ws.Range("E9").Value = 0.6230844082087783
$ws.Range("H9").Value = 1

$ws.Range("D10").Value = 0.4956050615108646
$ws.Range("E10").Value = 0.5043949384891353
$ws.Range("H10").Value = 1

$ws.Range("D11").Value = 0.1820489634123933
$ws.Range("E11").Value = 0.8179510365876066
$ws.Range("F11").Value = 0.6575940847396851
$ws.Range("H11").Value = 1

# Iterations = 200 block (rows 12-21): only H column is new, D/E unchanged
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1
